$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Naman"
$ws.Range("B6").Value = "Naman"
$ws.Range("C6").Value = $true

$ws.Range("D6").Select()
